$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update F2 to use a formula instead of the static value
$ws.Range("F2").Formula = "=5*5"

# Update the active selection to match the target state (F3)
$ws.Range("F3").Select()
